$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.930.61"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "2.647.69"
$ws.Range("E3").Value = "  +3.40%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'514.38"
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("D6").Value = "'144.03"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("E7").Value = "  -0.48%  "
$ws.Range("E8").Value = "  +1.78%  "
$ws.Range("D9").Value = "2.677.99"
$ws.Range("E9").Value = "  +4.60%  "
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("E11").Value = "  +3.50%  "
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("D14").Value = "3.113.40"
$ws.Range("E14").Value = "  +3.33%  "
$ws.Range("D15").Value = "58.931.42"
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("D16").Value = "'21.08"
$ws.Range("E16").Value = "  +2.54%  "
$ws.Range("E17").Value = "  +2.14%  "
$ws.Range("D18").Value = "2.674.52"
$ws.Range("E18").Value = "  +4.09%  "
$ws.Range("D19").Value = "'4.54"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").Value = "'340.47"
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("E21").Value = "  +3.38%  "
$ws.Range("E22").Value = "  +2.85%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'61.03"
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("D25").Value = "'0.419"
$ws.Range("E25").Value = "  +2.90%  "
$ws.Range("D26").Value = "2.759.31"
$ws.Range("E26").Value = "  +2.86%  "
$ws.Range("D27").Value = "'0.994"
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("E28").Value = "  +3.79%  "
$ws.Range("D29").Value = "0.0₃0808"
$ws.Range("E29").Value = "  +4.14%  "
$ws.Range("D30").Value = "'7.14"
$ws.Range("E30").Value = "  +4.29%  "
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").Value = "'6.38"
$ws.Range("E32").Value = "  +8.92%  "
$ws.Range("E33").Value = "  +2.66%  "
$ws.Range("D34").Value = "'18.91"
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("D35").Value = "'149.22"
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("E36").Value = "  +13.66%  "
$ws.Range("D37").Value = "'4.03"
$ws.Range("E37").Value = "  +4.45%  "
$ws.Range("E38").Value = "  +3.57%  "
$ws.Range("D39").Value = "'0.855"
$ws.Range("E39").Value = "  +4.49%  "
$ws.Range("D40").Value = "'36.62"
$ws.Range("E40").Value = "  +2.17%  "
$ws.Range("E41").Value = "  +4.28%  "
$ws.Range("D42").Value = "'1.40"
$ws.Range("E42").Value = "  +1.74%  "
$ws.Range("D43").Value = "'282.64"
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").Value = "'0.619"
$ws.Range("E44").Value = "  +2.07%  "
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("D46").Value = "'0.0981"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").Value = "'19.53"
$ws.Range("E47").Value = "  +5.02%  "
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("E50").Value = "  +4.83%  "
$ws.Range("D51").Value = "'10.27"
$ws.Range("E51").Value = "  -0.60%  "
